$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Discord ID column (D): convert text (scientific-notation) values to real numbers ---
# Row 2 keeps the default/general style (no explicit number format applied).
$ws.Range("D2").Value = 827240676667293056

# Rows 3-6 and 8 get an explicit scientific-notation number format.
$ws.Range("D3").NumberFormat = "0.00E+00"
$ws.Range("D3").Value = 832491214673543040

$ws.Range("D4").NumberFormat = "0.00E+00"
$ws.Range("D4").Value = 1277994753878789888

$ws.Range("D5").NumberFormat = "0.00E+00"
$ws.Range("D5").Value = 381439468764134976

$ws.Range("D6").NumberFormat = "0.00E+00"
$ws.Range("D6").Value = 837357129270034048

$ws.Range("D8").NumberFormat = "0.00E+00"
$ws.Range("D8").Value = 723479360379550976

# Row 9: Player ID (B) and Discord ID (D) become plain numbers (General format).
$ws.Range("B9").NumberFormat = "General"
$ws.Range("B9").Value = 137534547

$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Value = 173312357244600000

# --- Date Added column (H) for rows 7 & 8: replace static dates with formulas ---
$ws.Range("H7").Formula = '=TEXT(DATE(2025,4,25) + TIME(17,49,0), "yyyy-mm-dd hh:mm")'
$ws.Range("H8").Formula = '=TEXT(DATE(2025,4,24) + TIME(17,49,0), "yyyy-mm-dd hh:mm")'

# --- Update the active selection to match the saved workbook state ---
$ws.Range("B11").Select()
